# Applies the "Created Disc and Phenolic Spacer / Reduced length of airframe
# back to 310mm (removed space 2 x 2mm in inter-mount)" edit to Lengths.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New Notes (column D) text for the "Plywood disc" spacer component ---
$discNote = "Plywood disc component that has the Stepper motor mounted to"

# Row 4 (s2 mount) and Row 9 (s1 mount) both get the new "disc" note.
$ws.Range("D4").Value = $discNote
$ws.Range("D9").Value = $discNote

# --- Row 5 (s2 inter-mount spacer): was a formula 25+B4+B6 (=37), now a
#     fixed literal value of 35 (2mm shorter) ---
$ws.Range("B5").Value = 35

# Row 5 Notes -> new "Phenolic Tube that has s2 mount and s2 inter-mount
# installed in it" note, with "s2 mount" / "s2 inter-mount" in bold.
$s2Text = "Phenolic Tube that has s2 mount and s2 inter-mount installed in it"
$ws.Range("D5").Value = $s2Text
$ws.Range("D5").Characters(24, 8).Font.Bold = $true
$ws.Range("D5").Characters(32, 5).Font.Bold = $false
$ws.Range("D5").Characters(37, 14).Font.Bold = $true
$ws.Range("D5").Characters(51, 16).Font.Bold = $false

# --- Row 8 (s1 inter-mount spacer): was a formula 25+B7+B9 (=37), now a
#     fixed literal value of 35 (2mm shorter) ---
$ws.Range("B8").Value = 35

# Row 8 Notes -> new "Phenolic Tube that has s1 mount and s1 inter-mount
# installed in it" note, with "s1 mount" / "s1 inter-mount" in bold.
$s1Text = "Phenolic Tube that has s1 mount and s1 inter-mount installed in it"
$ws.Range("D8").Value = $s1Text
$ws.Range("D8").Characters(24, 8).Font.Bold = $true
$ws.Range("D8").Characters(32, 5).Font.Bold = $false
$ws.Range("D8").Characters(37, 14).Font.Bold = $true
$ws.Range("D8").Characters(51, 16).Font.Bold = $false

# --- Row 11 (interchange) Notes -> new phenolic-tube note ---
$ws.Range("D11").Value = "Phenolic tube that has interchange-lower at bottom and interchange-upper in it."

# --- Row 13 (tray-lower): 6 -> 9 ---
$ws.Range("B13").Value = 9

# --- Row 15 (tray-upper): 8 -> 9 ---
$ws.Range("B15").Value = 9

# --- Update selection / scroll position to match the saved view state ---
$ws.Range("C14").Select() | Out-Null

Write-Host "Lengths.xlsx updated"
